$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 12346864
$ws.Range("J19").Value = 15874329
$ws.Range("L19").Value = 15874329
$ws.Range("N19").Value = -15874679

$ws.Range("H34").Value = 13162.5
$ws.Range("I34").Value = 13162.5
$ws.Range("K34").Value = 13162.5
$ws.Range("M34").Value = -12959.5

$ws.Range("H36").Value = 13162.5
$ws.Range("I36").Value = 13162.5
$ws.Range("K36").Value = 13162.5
$ws.Range("M36").Value = -12447.5

$ws.Range("H53").Value = 18519714
$ws.Range("I53").Value = 66667680
$ws.Range("J53").Value = 1264.3077
$ws.Range("K53").Value = 66667680
$ws.Range("L53").Value = 1264.3077
$ws.Range("M53").Value = -66667043
$ws.Range("N53").Value = -2538.3077

$ws.Range("H98").Value = 1780.025
$ws.Range("I98").Value = 1119.742
$ws.Range("K98").Value = 1119.742
$ws.Range("M98").Value = 378.258

$ws.Range("H113").Value = 12500
$ws.Range("I113").Value = 20000
$ws.Range("K113").Value = 20000
$ws.Range("M113").Value = -16746

$ws.Range("H122").Value = 1780.025
$ws.Range("I122").Value = 1119.742
$ws.Range("K122").Value = 3359.226
$ws.Range("M122").Value = -909.2259999999997

$ws.Range("H129").Value = 2309.7693
$ws.Range("I129").Value = 1309.5
$ws.Range("J129").Value = 2609.85
$ws.Range("K129").Value = 3928.5
$ws.Range("L129").Value = 7829.549999999999
$ws.Range("M129").Value = 1071.5
$ws.Range("N129").Value = -17829.55

$ws.Range("H132").Value = 4312.524
$ws.Range("I132").Value = 4240.65
$ws.Range("K132").Value = 12721.95
$ws.Range("M132").Value = -10191.95

$ws.Range("H135").Value = 2055.8484
$ws.Range("I135").Value = 1962.3928
$ws.Range("K135").Value = 17661.5352
$ws.Range("M135").Value = -15126.5352

$ws.Range("H138").Value = 5288.12
$ws.Range("I138").Value = 2451.0625
$ws.Range("J138").Value = 6623.206
$ws.Range("K138").Value = 7353.1875
$ws.Range("L138").Value = 19869.618
$ws.Range("M138").Value = -2213.1875
$ws.Range("N138").Value = -30149.618

$ws.Range("H141").Value = 2675.0527
$ws.Range("I141").Value = 2379.2778
$ws.Range("K141").Value = 7137.8334
$ws.Range("M141").Value = -1957.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1151.4595
$ws.Range("I32").Value = 1200.9254
$ws.Range("K32").Value = 1200.9254
$ws.Range("M32").Value = -913.9254000000001

$ws.Range("H61").Value = 3212.3333
$ws.Range("I61").Value = 2913.4546
$ws.Range("K61").Value = 2913.4546
$ws.Range("M61").Value = -2701.4546

$ws.Range("H132").Value = 391745.47
$ws.Range("I132").Value = 635251.7
$ws.Range("J132").Value = 15417.637
$ws.Range("K132").Value = 1905755.1
$ws.Range("L132").Value = 46252.911
$ws.Range("M132").Value = -1903225.1
$ws.Range("N132").Value = -51312.911

$ws.Range("H136").Value = 3212.3333
$ws.Range("I136").Value = 2913.4546
$ws.Range("K136").Value = 8740.363799999999
$ws.Range("M136").Value = -6190.363799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 6000
$ws.Range("N86").Value = -8246
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 6000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 30000
$ws.Range("N89").Value = -41232
$ws.Range("M89").ClearContents()

$ws.Range("H134").Value = 40596.758
$ws.Range("I134").Value = 1837.0476
$ws.Range("J134").Value = 142341
$ws.Range("K134").Value = 5511.142800000001
$ws.Range("L134").Value = 427023
$ws.Range("M134").Value = -2976.142800000001
$ws.Range("N134").Value = -432093

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 973.46155
$ws.Range("I94").Value = 498.66666
$ws.Range("J94").Value = 1115.9
$ws.Range("K94").Value = 498.66666
$ws.Range("L94").Value = 1115.9
$ws.Range("M94").Value = -47.66665999999998
$ws.Range("N94").Value = -2017.9

$ws.Range("H132").Value = 3209.7144
$ws.Range("I132").Value = 1993.6
$ws.Range("K132").Value = 5980.799999999999
$ws.Range("M132").Value = -3450.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1770.1194
$ws.Range("J68").Value = 1828.4237
$ws.Range("L68").Value = 5485.2711
$ws.Range("N68").Value = -7107.2711

$ws.Range("H71").Value = 1770.1194
$ws.Range("J71").Value = 1828.4237
$ws.Range("L71").Value = 16455.8133
$ws.Range("N71").Value = -24567.8133

$ws.Range("H140").Value = 3554.4736
$ws.Range("I140").Value = 2248.75
$ws.Range("K140").Value = 6746.25
$ws.Range("M140").Value = -1566.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H132").Value = 59755.61
$ws.Range("I132").Value = 17589.477
$ws.Range("J132").Value = 502500
$ws.Range("K132").Value = 52768.431
$ws.Range("L132").Value = 1507500
$ws.Range("M132").Value = -50238.431
$ws.Range("N132").Value = -1512560

$ws.Range("H141").Value = 59900
$ws.Range("J141").Value = 59900
$ws.Range("L141").Value = 59900
$ws.Range("N141").Value = -70260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 726406.1
$ws.Range("I40").Value = 781668.1
$ws.Range("K40").Value = 781668.1
$ws.Range("M40").Value = -781532.1

$ws.Range("H46").Value = 2213.8572
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2213.8572
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2213.8572
$ws.Range("N46").Value = -2589.8572
$ws.Range("M46").ClearContents()

$ws.Range("H61").Value = 3250.205
$ws.Range("J61").Value = 4192.353
$ws.Range("L61").Value = 4192.353
$ws.Range("N61").Value = -4596.353

$ws.Range("H68").Value = 114832.445
$ws.Range("J68").Value = 203738.6
$ws.Range("L68").Value = 203738.6
$ws.Range("N68").Value = -205236.6

$ws.Range("H71").Value = 114832.445
$ws.Range("J71").Value = 203738.6
$ws.Range("L71").Value = 1018693
$ws.Range("N71").Value = -1026181

$ws.Range("H82").Value = 1499.2
$ws.Range("I82").Value = 422.57144
$ws.Range("J82").Value = 2441.25
$ws.Range("K82").Value = 422.57144
$ws.Range("L82").Value = 2441.25
$ws.Range("M82").Value = -61.57144
$ws.Range("N82").Value = -3163.25

$ws.Range("H85").Value = 1499.2
$ws.Range("I85").Value = 422.57144
$ws.Range("J85").Value = 2441.25
$ws.Range("K85").Value = 422.57144
$ws.Range("L85").Value = 2441.25
$ws.Range("M85").Value = 825.4285600000001
$ws.Range("N85").Value = -4937.25

$ws.Range("H100").Value = 231479.8
$ws.Range("J100").Value = 62500
$ws.Range("L100").Value = 62500
$ws.Range("N100").Value = -63582

$ws.Range("H113").Value = 3250.205
$ws.Range("J113").Value = 4192.353
$ws.Range("L113").Value = 4192.353
$ws.Range("N113").Value = -8532.352999999999

$ws.Range("H132").Value = 3468.9546
$ws.Range("I132").Value = 2559.3057
$ws.Range("K132").Value = 7677.9171
$ws.Range("M132").Value = -5147.9171

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 891.8
$ws.Range("J8").Value = 934.5
$ws.Range("L8").Value = 934.5
$ws.Range("N8").Value = -1214.5

$ws.Range("H104").Value = 55916.168
$ws.Range("J104").Value = 55916.168
$ws.Range("L104").Value = 55916.168
$ws.Range("N104").Value = -62904.168

$ws.Range("H126").Value = 11301
$ws.Range("I126").Value = 3079.8
$ws.Range("K126").Value = 9239.400000000001
$ws.Range("M126").Value = -6769.400000000001

$ws.Range("H132").Value = 17882.3
$ws.Range("I132").Value = 1212.2264
$ws.Range("K132").Value = 3636.6792
$ws.Range("M132").Value = -1106.6792
